$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the randomized sample values (D4:H7) ---
$ws.Range("D4").Value = -0.986
$ws.Range("E4").Value = 0.349
$ws.Range("F4").Value = -0.377
$ws.Range("G4").Value = -0.165
$ws.Range("H4").Value = -0.402

$ws.Range("D5").Value = 0.48
$ws.Range("E5").Value = 0.23
$ws.Range("F5").Value = -0.28
$ws.Range("G5").Value = 0.404
$ws.Range("H5").Value = 0.037

$ws.Range("D6").Value = 5.992
$ws.Range("E6").Value = 4.711
$ws.Range("F6").Value = 2.504
$ws.Range("G6").Value = 0.04
$ws.Range("H6").Value = -0.211

$ws.Range("D7").Value = 6.521
$ws.Range("E7").Value = 3.574
$ws.Range("F7").Value = 1.757
$ws.Range("G7").Value = 1.482
$ws.Range("H7").Value = 0.081

# --- Apply a 3-decimal numeric format to the updated block ---
$ws.Range("D4:H7").NumberFormat = "0.000"

# --- Move the active selection ---
$ws.Range("G17").Select()
